# Apply the edits described in the commit:
# "BOM added, some calculations added in XLSX"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared string "Подбор" referenced from D6 and D16 (Bill of Materials item)
$ws.Range("D6").Value = "Подбор"
$ws.Range("D16").Value = "Подбор"

# Update baud/bit values in the second timing block
$ws.Range("B15").Value = 7
$ws.Range("B16").Value = 25

# B21 no longer mirrors B19 via formula; it becomes a fixed experimental value
$ws.Range("B21").Value = 6.8

# Update the view: scrolled down with new active selection cell B24
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("B24").Select()
